$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) figures for each
# coin row (rows 2-51) with the latest scraped values.
# Column D is forced to Text format ("@") before assignment so that values
# such as "1.00", "0.517" or "239.77" are preserved exactly as literal text
# instead of being auto-converted into numbers by Excel (which would drop
# trailing zeros / introduce floating point noise). Column E values already
# contain padding spaces and a percent sign so Excel keeps them as text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.167.43'
$ws.Range("E2").Value = '  +2.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.300.98'
$ws.Range("E3").Value = '  +2.37%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.75'
$ws.Range("E5").Value = '  +0.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.94'
$ws.Range("E6").Value = '  +7.63%  '

$ws.Range("E7").Value = '  +0.54%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +6.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.67'
$ws.Range("E10").Value = '  +2.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.53'
$ws.Range("E11").Value = '  +1.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0808'
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("E13").Value = '  -1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.92'
$ws.Range("E14").Value = '  +2.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.659.68'
$ws.Range("E15").Value = '  +2.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.07'
$ws.Range("E16").Value = '  +4.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.307.75'
$ws.Range("E17").Value = '  +1.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.797'
$ws.Range("E18").Value = '  +2.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.134.60'
$ws.Range("E19").Value = '  +2.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.94'
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0919'
$ws.Range("E21").Value = '  +2.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  +4.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.66'
$ws.Range("E23").Value = '  +0.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '239.77'
$ws.Range("E24").Value = '  +1.83%  '

$ws.Range("E25").Value = '  +3.17%  '

$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.64'
$ws.Range("E28").Value = '  +5.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.33'
$ws.Range("E29").Value = '  +10.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.92'
$ws.Range("E30").Value = '  -5.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.53'
$ws.Range("E31").Value = '  +0.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.69'
$ws.Range("E32").Value = '  -2.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.21'
$ws.Range("E33").Value = '  +1.24%  '

$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.19'
$ws.Range("E35").Value = '  +4.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.53'
$ws.Range("E36").Value = '  +6.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0729'
$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.99'
$ws.Range("E38").Value = '  -2.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.48'
$ws.Range("E39").Value = '  +10.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  +2.91%  '

$ws.Range("E41").Value = '  +3.57%  '

$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.53'
$ws.Range("E43").Value = '  +15.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0288'
$ws.Range("E44").Value = '  +2.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.958.41'
$ws.Range("E45").Value = '  +0.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.67'
$ws.Range("E46").Value = '  +1.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.05'
$ws.Range("E47").Value = '  +5.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.15'
$ws.Range("E48").Value = '  +5.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '57.01'
$ws.Range("E49").Value = '  +6.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.91'
$ws.Range("E50").Value = '  +0.50%  '

$ws.Range("E51").Value = '  +8.21%  '

